$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D: sex = "U" for rows 2-11
$ws.Range("D2:D11").Value = "U"

# Column F: sire values for rows 2-11
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 3
$ws.Range("F4").Value = 5
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 5
$ws.Range("F7").Value = 6
$ws.Range("F8").Value = 8
$ws.Range("F9").Value = 10
$ws.Range("F10").Value = 8
$ws.Range("F11").Value = 10

# Column G: dam values for rows 2-11
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 5
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 6
$ws.Range("G8").Value = 6
$ws.Range("G9").Value = 6
$ws.Range("G10").Value = 10
$ws.Range("G11").Value = 10

# Update the selection to reflect the new used range
$ws.Range("D2:G11").Select()

Write-Output "done"
